# Insert a new data row before row 49 (pushing the existing rows 49-97
# down to 50-98, and growing the sheet's used range to A1:R98), then
# populate the newly inserted row 49 with its own record.
#
# This mirrors: one new weekly "Cilantro" price quote was added to the
# top of the historical series kept by this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 49:97 -> 50:98, leaving a blank row 49 behind.
$ws.Rows("49:49").Insert()

# Fill in the new row 49.
$ws.Cells.Item(49, 1).Value  = 1
$ws.Cells.Item(49, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(49, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(49, 4).Value  = 44907
$ws.Cells.Item(49, 5).Value  = 15
$ws.Cells.Item(49, 6).Value  = 100112040
$ws.Cells.Item(49, 7).Value  = "Cilantro"
$ws.Cells.Item(49, 8).Value  = "Sin especificar"
$ws.Cells.Item(49, 9).Value  = "Primera"
$ws.Cells.Item(49, 10).Value = 500
$ws.Cells.Item(49, 11).Value = 1300
$ws.Cells.Item(49, 12).Value = 1500
$ws.Cells.Item(49, 13).Value = 1400
$ws.Cells.Item(49, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(49, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(49, 16).Value = 700
$ws.Cells.Item(49, 17).Value = 2
$ws.Cells.Item(49, 18).Value = "Hortaliza"
